$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.195.82"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.73%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.438.08"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.51%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +1.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9159"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -7.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "275.11"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.67%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3620"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.76%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3071"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.91%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "38.66"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.96%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.021"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.08%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06481"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.52%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9996"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.38%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.322"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.39"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.028"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001005"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.438.32"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.78%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9321"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -6.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05616"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.38"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -6.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.389"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.60%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.20"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.59%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.82"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.242"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.53%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.213.48"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.70%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "137.74"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.94%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.110"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -7.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.82"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.86%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.591.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "109.71"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.933"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -5.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8014"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.808"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -10.17%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07634"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.468"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.75%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05830"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.645"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.126"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01977"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.13"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.53%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1842"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9261"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -6.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.076"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -15.88%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5188"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.26%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.485"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.86"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "116.32"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.94%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5059"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.14%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.724"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.94%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06381"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9888"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.31%  "
